$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up the "D" (column E) imputed values that were left blank / wrongly filled
$ws.Range("E3").Value = -5.7
$ws.Range("E5").Value = ""
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("E23").Value = ""

# Row 26 ("RM 232") was a bad / duplicate row - remove it entirely, shifting
# everything below up by one.
$ws.Rows(26).Delete()

# After that shift, the row that used to be "SC 92" (old row 28) is now row 27
# - remove it too.
$ws.Rows(27).Delete()

# The row that is now "SC 193" (old row 34, now row 32) was missing its D
# value - fill it in.
$ws.Range("E32").Value = -6.4
